$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be read as text before assigning values, so that
# numeric-looking price strings (e.g. "318.30") are not auto-converted
# to numbers by Excel; this mirrors the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "44.117.62"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "2.276.12"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "318.30"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "103.32"
$ws.Range("E6").Value = "  +5.93%  "
$ws.Range("E7").Value = "  +1.70%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("D10").Value = "38.72"
$ws.Range("E10").Value = "  +6.42%  "
$ws.Range("D11").Value = "0.0840"
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("D12").Value = "7.87"
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "2.622.73"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").Value = "0.875"
$ws.Range("D16").Value = "14.55"
$ws.Range("E16").Value = "  +3.87%  "
$ws.Range("D17").Value = "2.273.45"
$ws.Range("D18").Value = "44.006.84"
$ws.Range("E18").Value = "  +2.92%  "
$ws.Range("D19").Value = "14.34"
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("D20").Value = "0.0₃0999"
$ws.Range("E20").Value = "  +4.27%  "
$ws.Range("D21").Value = "6.66"
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("D22").Value = "66.18"
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").Value = "237.93"
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("E25").Value = "  +4.00%  "
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").Value = "39.41"
$ws.Range("E28").Value = "  +17.35%  "
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").Value = "6.52"
$ws.Range("E30").Value = "  +4.46%  "
$ws.Range("D31").Value = "162.03"
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("D32").Value = "20.50"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").Value = "0.0879"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "3.27"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").Value = "4.52"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("E39").Value = "  +4.96%  "
$ws.Range("D40").Value = "3.86"
$ws.Range("E40").Value = "  +4.89%  "
$ws.Range("D41").Value = "15.60"
$ws.Range("E41").Value = "  +29.06%  "
$ws.Range("D42").Value = "0.0327"
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "1.778.92"
$ws.Range("E44").Value = "  -5.33%  "
$ws.Range("D45").Value = "0.208"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").Value = "5.42"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("D47").Value = "84.94"
$ws.Range("E47").Value = "  -3.55%  "
$ws.Range("D48").Value = "8.91"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D49").Value = "59.55"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").Value = "74.68"
$ws.Range("E50").Value = "  -4.21%  "
$ws.Range("D51").Value = "104.46"
$ws.Range("E51").Value = "  +3.73%  "

# Restore default (no explicit) style on column D so the saved XML
# matches the original cells, which carried no style attribute.
$ws.Range("D2:D51").Style = "Normal"
